# Updated cryptos list (Price / Volume(1h) columns) with refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "27.886.54" },
    @{ Cell = "E2"; Value = "  -4.47%  " },
    @{ Cell = "D3"; Value = "1.737.50" },
    @{ Cell = "E3"; Value = "  -4.84%  " },
    @{ Cell = "D4"; Value = "1.002" },
    @{ Cell = "E4"; Value = "  -0.20%  " },
    @{ Cell = "D5"; Value = "226.51" },
    @{ Cell = "E5"; Value = "  -3.48%  " },
    @{ Cell = "D6"; Value = "0.5768" },
    @{ Cell = "E6"; Value = "  -3.82%  " },
    @{ Cell = "E7"; Value = "  -0.02%  " },
    @{ Cell = "D8"; Value = "0.2734" },
    @{ Cell = "E8"; Value = "  -0.94%  " },
    @{ Cell = "D9"; Value = "23.07" },
    @{ Cell = "E9"; Value = "  -1.31%  " },
    @{ Cell = "D10"; Value = "0.06622" },
    @{ Cell = "E10"; Value = "  -4.32%  " },
    @{ Cell = "D11"; Value = "0.07545" },
    @{ Cell = "E11"; Value = "  -0.84%  " },
    @{ Cell = "D12"; Value = "1.749.00" },
    @{ Cell = "E12"; Value = "  -4.81%  " },
    @{ Cell = "D13"; Value = "4.704" },
    @{ Cell = "E13"; Value = "  -0.25%  " },
    @{ Cell = "D14"; Value = "0.6017" },
    @{ Cell = "E14"; Value = "  -3.52%  " },
    @{ Cell = "D15"; Value = "1.975.34" },
    @{ Cell = "E15"; Value = "  -4.71%  " },
    @{ Cell = "D16"; Value = "74.45" },
    @{ Cell = "E16"; Value = "  -3.52%  " },
    @{ Cell = "D17"; Value = "0.000008665" },
    @{ Cell = "E17"; Value = "  -11.34%  " },
    @{ Cell = "D18"; Value = "27.873.35" },
    @{ Cell = "E19"; Value = "  -4.01%  " },
    @{ Cell = "E20"; Value = "  -0.16%  " },
    @{ Cell = "D21"; Value = "204.84" },
    @{ Cell = "E21"; Value = "  -4.79%  " },
    @{ Cell = "D22"; Value = "11.27" },
    @{ Cell = "E22"; Value = "  -2.32%  " },
    @{ Cell = "D23"; Value = "6.618" },
    @{ Cell = "E23"; Value = "  -3.07%  " },
    @{ Cell = "E24"; Value = "  -0.03%  " },
    @{ Cell = "D25"; Value = "150.13" },
    @{ Cell = "E25"; Value = "  -3.62%  " },
    @{ Cell = "D26"; Value = "8.021" },
    @{ Cell = "E26"; Value = "  +1.01%  " },
    @{ Cell = "D27"; Value = "0.1232" },
    @{ Cell = "E27"; Value = "  -4.28%  " },
    @{ Cell = "D28"; Value = "16.16" },
    @{ Cell = "E28"; Value = "  -1.78%  " },
    @{ Cell = "D29"; Value = "0.06195" },
    @{ Cell = "E29"; Value = "  -4.67%  " },
    @{ Cell = "D30"; Value = "1.381" },
    @{ Cell = "E30"; Value = "  -3.24%  " },
    @{ Cell = "D31"; Value = "1.393" },
    @{ Cell = "E31"; Value = "  -3.32%  " },
    @{ Cell = "D32"; Value = "3.739" },
    @{ Cell = "E32"; Value = "  -1.72%  " },
    @{ Cell = "D33"; Value = "3.727" },
    @{ Cell = "E33"; Value = "  -1.21%  " },
    @{ Cell = "D34"; Value = "1.679" },
    @{ Cell = "E34"; Value = "  -2.12%  " },
    @{ Cell = "D35"; Value = "1.035" },
    @{ Cell = "E35"; Value = "  -4.90%  " },
    @{ Cell = "D36"; Value = "0.6388" },
    @{ Cell = "E36"; Value = "  -0.57%  " },
    @{ Cell = "D37"; Value = "2.442" },
    @{ Cell = "E37"; Value = "  -3.97%  " },
    @{ Cell = "D38"; Value = "2.711" },
    @{ Cell = "E38"; Value = "  -1.85%  " },
    @{ Cell = "E39"; Value = "  -4.45%  " },
    @{ Cell = "D40"; Value = "1.119.86" },
    @{ Cell = "E40"; Value = "  -1.03%  " },
    @{ Cell = "D41"; Value = "6.167" },
    @{ Cell = "E41"; Value = "  -6.31%  " },
    @{ Cell = "D42"; Value = "0.8733" },
    @{ Cell = "E42"; Value = "  -1.53%  " },
    @{ Cell = "E43"; Value = "  +0.17%  " },
    @{ Cell = "D44"; Value = "100.08" },
    @{ Cell = "E44"; Value = "  -0.52%  " },
    @{ Cell = "D45"; Value = "1.886.46" },
    @{ Cell = "E45"; Value = "  -4.89%  " },
    @{ Cell = "D46"; Value = "59.30" },
    @{ Cell = "E46"; Value = "  -4.13%  " },
    @{ Cell = "D47"; Value = "0.00000000110" },
    @{ Cell = "E47"; Value = "  -3.27%  " },
    @{ Cell = "D48"; Value = "1.572" },
    @{ Cell = "E48"; Value = "  -2.09%  " },
    @{ Cell = "D49"; Value = "8.281" },
    @{ Cell = "E49"; Value = "  -2.09%  " },
    @{ Cell = "D50"; Value = "0.05376" },
    @{ Cell = "E50"; Value = "  -2.30%  " },
    @{ Cell = "D51"; Value = "0.4415" },
    @{ Cell = "E51"; Value = "  -2.53%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}